$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update progress percentages for the "MOVIMENTAÇÃO" / Bruno row (row 11)
# under the LIST / CREATE / SHOW / EDIT columns (E:H)
$ws.Range("E11").Value = 80
$ws.Range("F11").Value = 90
$ws.Range("G11").Value = 90
$ws.Range("H11").Value = 90

# Reflect the user's last active selection, as recorded in the saved file
$ws.Range("E13").Select()
